$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 42
$ws.Range("I2").Value = 92
$ws.Range("J2").Value = 471
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 148
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 89
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = 84
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 738
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 740
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 9
$ws.Range("AA2").Value = 0
